$wb = $excel.ActiveWorkbook

# Both "展览" and "全部类型" sheets contain the same 4 event rows (rows 2-5)
# with an updated "想去人数" (interested-count) column F.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 39
    $ws.Range("F3").Value = 80
    $ws.Range("F4").Value = 2090
    $ws.Range("F5").Value = 166
}
